$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily row appended below the existing data (A1:D62 -> A1:D63).
# Date/weekday are stored as literal text (matching the rest of the
# sheet, which never holds real Excel date serials) — force text entry
# via NumberFormat, then drop back to the default "Normal" style so no
# stray formatting is left on the new cell.
$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "2025/10/05"
$ws.Range("A63").Style = "Normal"

$ws.Range("B63").Value = "日"
$ws.Range("C63").Value = 8
$ws.Range("D63").Value = 45
